$d = $word.ActiveDocument

# The document starts with two empty paragraphs: the first one plain,
# the second one styled "Normal". The edit types "Entrevista zero" into
# the first paragraph (leaving Word's automatic "last edit" _GoBack
# bookmark collapsed right after the new text) and clears the explicit
# "Normal" style from the now-trailing empty paragraph.

$first = $d.Paragraphs(1).Range
$first.Collapse(1)
$d.Bookmarks.Add("_GoBack", $first)
$first.InsertBefore("Entrevista zero")

$d.Paragraphs(2).Range.ParagraphFormat.Style = "Normal"
